# 07April2022 Selenium DataDriven part2
#
# - Rename "Sheet1" -> "sheet1"
# - Add two new trailing sheets: "Sheet3" (blank) and "Sheet4" (one cell: "dsds")
# - TestData sheet:
#     * B6 gains an extra column of data (col S already existed) -> autoFilter/_FilterDatabase
#       widened from A1:R8 to A1:S8
#     * C6 "Hemanth" -> "Geetha"
#     * H6 "Male" -> "Female"
#     * Active selection on TestData moves to H6
# - "sheet1" sheet gains a selection at E26
# - Defined name TestData!_FilterDatabase widened to $A$1:$S$8

$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> sheet1 -----------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Name = "sheet1"

$testData = $wb.Worksheets.Item("TestData")

# --- Add two new worksheets at the end of the workbook ---------------------
# (Sheet4!A1 = "dsds" must land in the shared-string table BEFORE TestData!C6
#  is changed to "Geetha", so set it first.)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sheet3.Name = "Sheet3"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$sheet4.Name = "Sheet4"
$sheet4.Range("A1").Value = "dsds"

# --- TestData edits ----------------------------------------------------
$testData.Range("C6").Value = "Geetha"
$testData.Range("H6").Value = "Female"

# Re-apply the AutoFilter so it covers the full A1:S8 range (it already did,
# this just keeps the filter element / mode consistent after the edits).
if ($testData.AutoFilterMode) {
    $testData.AutoFilterMode = $false
}
$testData.Range("A1:S8").AutoFilter()

# Widen the hidden _FilterDatabase defined name to match.
$testData.Names.Item(1).RefersTo = "=TestData!`$A`$1:`$S`$8"

# --- Selections --------------------------------------------------------
# sheet1: select E26
$ws1.Activate()
$ws1.Range("E26").Select()

# TestData: select H6, and leave it the active tab (matches original state).
$testData.Activate()
$testData.Range("H6").Select()
